$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $r = $ws.Range($Cell)
    $r.NumberFormat = "@"
    $r.Value = $Text
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "27.782.61"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "1.649.00"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue "D5" "213.59"
$ws.Range("E5").Value = "  +0.15%  "
Set-TextValue "D6" "0.533"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("E9").Value = "  -0.77%  "
Set-TextValue "D11" "0.0891"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "1.881.95"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "1.648.30"
$ws.Range("E13").Value = "  -0.32%  "
Set-TextValue "D14" "4.05"
$ws.Range("E14").Value = "  -0.61%  "
Set-TextValue "D15" "0.565"
$ws.Range("E15").Value = "  -0.61%  "
Set-TextValue "D16" "64.46"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D17").Value = "27.750.40"
$ws.Range("E17").Value = "  +1.33%  "
Set-TextValue "D18" "232.61"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D19" "7.71"
$ws.Range("E19").Value = "  +3.41%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0726"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  -0.47%  "
Set-TextValue "D23" "10.16"
$ws.Range("E23").Value = "  +9.61%  "
$ws.Range("E24").Value = "  -3.51%  "
Set-TextValue "D25" "150.51"
$ws.Range("E25").Value = "  +2.31%  "
Set-TextValue "D26" "6.99"
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D33" "3.17"
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.440.25"
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("E37").Value = "  +0.62%  "
Set-TextValue "D38" "0.889"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("E39").Value = "  -0.71%  "
Set-TextValue "D40" "0.884"
$ws.Range("E40").Value = "  +12.12%  "
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D42" "1.00"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "5.63"
$ws.Range("E43").Value = "  +1.32%  "
Set-TextValue "D44" "66.58"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("D47").Value = "1.791.29"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("E48").Value = "  +3.34%  "
Set-TextValue "D49" "86.55"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("E50").Value = "  +2.38%  "
Set-TextValue "D51" "0.1000"
$ws.Range("E51").Value = "  -1.47%  "
